$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values are assigned with a leading apostrophe so Excel stores them as literal
# text (matching the source data, which includes multi-dot numbers, leading
# zeros, and percentage strings) without altering the number format of the cell.

$ws.Range('D2').Value = "'59.040.75"
$ws.Range('E2').Value = "'  +1.53%  "
$ws.Range('D3').Value = "'2.586.66"
$ws.Range('E3').Value = "'  -0.25%  "
$ws.Range('E4').Value = "'  -0.19%  "
$ws.Range('D5').Value = "'528.46"
$ws.Range('E5').Value = "'  +1.49%  "
$ws.Range('D6').Value = "'139.07"
$ws.Range('E6').Value = "'  -2.70%  "
$ws.Range('E7').Value = "'  +0.10%  "
$ws.Range('D8').Value = "'0.565"
$ws.Range('E8').Value = "'  -0.22%  "
$ws.Range('D9').Value = "'2.599.58"
$ws.Range('E9').Value = "'  -0.49%  "
$ws.Range('E10').Value = "'  -1.20%  "
$ws.Range('E11').Value = "'  -0.09%  "
$ws.Range('E12').Value = "'  -2.25%  "
$ws.Range('E13').Value = "'  +3.06%  "
$ws.Range('D14').Value = "'3.045.48"
$ws.Range('E14').Value = "'  -0.18%  "
$ws.Range('D15').Value = "'58.990.49"
$ws.Range('E15').Value = "'  +1.24%  "
$ws.Range('D16').Value = "'20.47"
$ws.Range('E16').Value = "'  +0.61%  "
$ws.Range('B17').Value = "'WrappedEther"
$ws.Range('C17').Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('D17').Value = "'2.604.25"
$ws.Range('E17').Value = "'  -0.40%  "
$ws.Range('B18').Value = "'ShibaInu"
$ws.Range('C18').Value = "'https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range('D18').Value = "'0.0000132"
$ws.Range('E18').Value = "'  -0.84%  "
$ws.Range('D19').Value = "'344.76"
$ws.Range('E19').Value = "'  +1.78%  "
$ws.Range('E20').Value = "'  -0.41%  "
$ws.Range('E21').Value = "'  -1.33%  "
$ws.Range('D22').Value = "'6.42"
$ws.Range('E22').Value = "'  -0.49%  "
$ws.Range('E23').Value = "'  +0.05%  "
$ws.Range('D24').Value = "'67.23"
$ws.Range('E24').Value = "'  +3.01%  "
$ws.Range('E25').Value = "'  -0.24%  "
$ws.Range('D26').Value = "'0.404"
$ws.Range('E26').Value = "'  +0.22%  "
$ws.Range('E27').Value = "'  +0.12%  "
$ws.Range('D28').Value = "'7.05"
$ws.Range('E28').Value = "'  +0.15%  "
$ws.Range('E29').Value = "'  +0.10%  "
$ws.Range('D30').Value = "'0.0₃0716"
$ws.Range('E30').Value = "'  -3.27%  "
$ws.Range('E31').Value = "'  +1.10%  "
$ws.Range('D32').Value = "'5.87"
$ws.Range('E32').Value = "'  -3.91%  "
$ws.Range('E33').Value = "'  -0.21%  "
$ws.Range('D34').Value = "'148.95"
$ws.Range('E34').Value = "'  -0.29%  "
$ws.Range('D35').Value = "'3.95"
$ws.Range('E35').Value = "'  -1.13%  "
$ws.Range('D36').Value = "'1.11"
$ws.Range('D37').Value = "'36.66"
$ws.Range('E37').Value = "'  +1.67%  "
$ws.Range('D38').Value = "'1.47"
$ws.Range('E38').Value = "'  +0.68%  "
$ws.Range('D39').Value = "'0.821"
$ws.Range('E39').Value = "'  -3.77%  "
$ws.Range('D40').Value = "'0.808"
$ws.Range('E40').Value = "'  -5.36%  "
$ws.Range('D41').Value = "'3.51"
$ws.Range('E41').Value = "'  -0.31%  "
$ws.Range('E42').Value = "'  +0.24%  "
$ws.Range('E43').Value = "'  -1.26%  "
$ws.Range('D44').Value = "'10.76"
$ws.Range('E44').Value = "'  +0.74%  "
$ws.Range('D45').Value = "'267.32"
$ws.Range('E45').Value = "'  -1.23%  "
$ws.Range('E46').Value = "'  -0.26%  "
$ws.Range('D47').Value = "'0.0513"
$ws.Range('E47').Value = "'  -1.40%  "
$ws.Range('D48').Value = "'18.29"
$ws.Range('E48').Value = "'  -2.26%  "
$ws.Range('D49').Value = "'1.957.17"
$ws.Range('E49').Value = "'  -0.44%  "
$ws.Range('E50').Value = "'  -0.68%  "
$ws.Range('D51').Value = "'18.13"
$ws.Range('E51').Value = "'  -0.52%  "
